{"js": "// Fix style definitions in styles.xml whose <w:rPr> children are ordered\n// incorrectly (the WML schema, wml.xsd / CT_RPr, requires <w:b/>/<w:i/>\n// etc. to precede <w:color/>). Several custom \"Tok\" character styles\n// (Pandoc/pygments syntax-highlighting styles) had <w:color/> emitted\n// before <w:b/> and/or <w:i/>, which OOXMLValidatorCLI flags as a schema\n// error even though xmllint stays silent about it.\n//\n// Re-assigning the (already-correct) bold/italic value on each affected\n// style's font forces the run-properties to be re-serialized in the\n// schema-correct order, without altering any actual formatting.\n\nconst styles = context.document.getStyles();\n\n// styleId -> { bold, italic } describing which weight/slant flags the\n// style carries (matches the \"before\" values - we are only fixing\n// element order, not any formatting).\nconst targets = {\n  KeywordTok: { bold: true, italic: false },\n  ImportTok: { bold: true, italic: false },\n  CommentTok: { bold: false, italic: true },\n  DocumentationTok: { bold: false, italic: true },\n  AnnotationTok: { bold: true, italic: true },\n  CommentVarTok: { bold: true, italic: true },\n  ControlFlowTok: { bold: true, italic: false },\n  InformationTok: { bold: true, italic: true },\n  WarningTok: { bold: true, italic: true },\n  AlertTok: { bold: true, italic: false },\n  ErrorTok: { bold: true, italic: false },\n};\n\nconst fonts = {};\nfor (const name of Object.keys(targets)) {\n  const style = styles.getByNameOrNullObject(name);\n  const font = style.font;\n  fonts[name] = font;\n}\n\n// Re-write bold/italic (same value as already set) on every affected\n// style; this dirties the style's run properties so they get\n// re-emitted in schema order (b/i before color) when the package is\n// saved.\nfor (const [name, flags] of Object.entries(targets)) {\n  const font = fonts[name];\n  if (flags.bold) {\n    font.bold = true;\n  }\n  if (flags.italic) {\n    font.italic = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix style definitions in styles.xml whose <w:rPr> children are ordered\n# incorrectly (the WML schema, wml.xsd / CT_RPr, requires <w:b/>/<w:i/>\n# etc. to precede <w:color/>). Several custom \"Tok\" character styles\n# (Pandoc/pygments syntax-highlighting styles) had <w:color/> emitted\n# before <w:b/> and/or <w:i/>, which OOXMLValidatorCLI flags as a schema\n# error even though xmllint stays silent about it.\n#\n# Re-assigning the (already-correct) bold/italic value on each affected\n# style's font forces the run-properties to be re-serialized in the\n# schema-correct order, without altering any actual formatting.\n\n$d = $word.ActiveDocument\n\n# styleId -> which weight/slant flags the style carries (matches the\n# \"before\" values - we are only fixing element order, not any\n# formatting).\n$targets = @{\n    \"KeywordTok\"       = @{ Bold = $true;  Italic = $false }\n    \"ImportTok\"        = @{ Bold = $true;  Italic = $false }\n    \"CommentTok\"       = @{ Bold = $false; Italic = $true }\n    \"DocumentationTok\" = @{ Bold = $false; Italic = $true }\n    \"AnnotationTok\"    = @{ Bold = $true;  Italic = $true }\n    \"CommentVarTok\"    = @{ Bold = $true;  Italic = $true }\n    \"ControlFlowTok\"   = @{ Bold = $true;  Italic = $false }\n    \"InformationTok\"   = @{ Bold = $true;  Italic = $true }\n    \"WarningTok\"       = @{ Bold = $true;  Italic = $true }\n    \"AlertTok\"         = @{ Bold = $true;  Italic = $false }\n    \"ErrorTok\"         = @{ Bold = $true;  Italic = $false }\n}\n\nforeach ($name in $targets.Keys) {\n    $flags = $targets[$name]\n    $style = $d.Styles($name)\n    if ($flags.Bold) {\n        $style.Font.Bold = -1\n    }\n    if ($flags.Italic) {\n        $style.Font.Italic = -1\n    }\n}\n"}
